$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift old row 5 ("社内PC..." posting) down to row 9 by inserting 4 rows at 5..8
$ws.Range("A5:A8").EntireRow.Insert()

# Drop the (now stale/incomplete) hyperlinks collection; rebuilt fresh below
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-10-31 18:25:45'
$ws.Range("B2").Value = '【急募】映像解析AIによる自動検出・モザイク処理スクリプト開発'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5424032'
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

# Row 3
$ws.Range("A3").Value = '2025-10-31 18:25:45'
$ws.Range("B3").Value = '【急募】Google Workspace/LLM連携!AI診断レポート自動生成システム構築'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5424558'
$ws.Range("G3").Value = 325
$ws.Range("H3").Value = '🔥AI,Ai'

# Row 4
$ws.Range("A4").Value = '2025-10-31 18:25:45'
$ws.Range("B4").Value = '【急募】ebayAPIを活用したShippingポリシー設定の専門家募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5415908'
$ws.Range("G4").Value = 183
$ws.Range("H4").Value = '🔥API'

# Row 5
$ws.Range("A5").Value = '2025-10-31 18:25:45'
$ws.Range("B5").Value = '【急募】Shopifyレンタルシステムのバックエンド開発者募集'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5424578'
$ws.Range("G5").Value = 83
$ws.Range("H5").Value = '◆開発'

# Row 6
$ws.Range("A6").Value = '2025-10-31 18:25:45'
$ws.Range("B6").Value = '【急募】Accessシステム改修・CSV読込・MySQLクラウド化【出張希望】'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5422936'
$ws.Range("G6").Value = 48
$ws.Range("H6").Value = '◇MySQL'

# Row 7
$ws.Range("A7").Value = '2025-10-31 18:25:45'
$ws.Range("B7").Value = 'eBayテラピークでのキーワード検索結果等の取得するためのシステム制作'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5420779'
$ws.Range("G7").Value = 33

# Row 8
$ws.Range("A8").Value = '2025-10-31 18:25:45'
$ws.Range("B8").Value = '【急募】クラファン制作代行であなたのアイデアを形に!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5424826'
$ws.Range("G8").Value = 18

# Row 9
$ws.Range("A9").Value = '2025-10-31 18:25:45'
$ws.Range("B9").Value = '【急募】社内PCのデータフォレンジック業務をお任せします'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5424258'
$ws.Range("G9").Value = 18

# Rebuild hyperlinks F2:F9 in order, reusing the shared "Hyperlink" cell style (s=1)
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5424032') | Out-Null
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5424558') | Out-Null
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5415908') | Out-Null
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5424578') | Out-Null
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5422936') | Out-Null
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5420779') | Out-Null
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5424826') | Out-Null
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5424258') | Out-Null
$ws.Range("F9").Style = "Hyperlink"
